# Refresh the cryptos price/volume list (GitHub Actions daily scrape).
# Price cells in column D that look like plain numbers are written with a
# leading apostrophe so Excel keeps them as text (preserving trailing
# zeros / thousands-dot formatting) instead of coercing them to floats.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '57.681.41'
$ws.Range('E2').Value = '  -0.45%  '
$ws.Range('D3').Value = '3.094.83'
$ws.Range('E3').Value = '  +0.62%  '
$ws.Range('E4').Value = '  +0.10%  '
$ws.Range('D5').Value = '''524.93'
$ws.Range('E5').Value = '  +1.78%  '
$ws.Range('D6').Value = '''141.06'
$ws.Range('E6').Value = '  -0.22%  '
$ws.Range('E7').Value = '  +0.10%  '
$ws.Range('D8').Value = '3.093.32'
$ws.Range('E8').Value = '  +0.60%  '
$ws.Range('D9').Value = '''0.440'
$ws.Range('E9').Value = '  +1.16%  '
$ws.Range('D10').Value = '''7.19'
$ws.Range('E10').Value = '  -1.76%  '
$ws.Range('E11').Value = '  +0.43%  '
$ws.Range('D12').Value = '''0.390'
$ws.Range('E12').Value = '  +3.26%  '
$ws.Range('D13').Value = '3.637.96'
$ws.Range('E13').Value = '  +0.99%  '
$ws.Range('E14').Value = '  +1.95%  '
$ws.Range('D15').Value = '''25.44'
$ws.Range('E15').Value = '  -4.45%  '
$ws.Range('D16').Value = '''0.0000163'
$ws.Range('E16').Value = '  -0.14%  '
$ws.Range('D17').Value = '57.914.63'
$ws.Range('E17').Value = '  -0.06%  '
$ws.Range('D18').Value = '3.086.77'
$ws.Range('E18').Value = '  +0.53%  '
$ws.Range('D19').Value = '''6.08'
$ws.Range('E19').Value = '  -1.31%  '
$ws.Range('D20').Value = '''12.74'
$ws.Range('E20').Value = '  -0.79%  '
$ws.Range('D21').Value = '''7.97'
$ws.Range('E21').Value = '  -1.40%  '
$ws.Range('D22').Value = '''340.47'
$ws.Range('E22').Value = '  +1.97%  '
$ws.Range('D23').Value = '''1.00'
$ws.Range('E23').Value = '  -0.01%  '
$ws.Range('D24').Value = '''0.511'
$ws.Range('E24').Value = '  +1.82%  '
$ws.Range('D25').Value = '''67.21'
$ws.Range('E25').Value = '  +3.32%  '
$ws.Range('E26').Value = '  -1.73%  '
$ws.Range('E27').Value = '  +0.19%  '
$ws.Range('D28').Value = '0.0₃0913'
$ws.Range('E28').Value = '  +0.59%  '
$ws.Range('B29').Value = 'USDe'
$ws.Range('C29').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D29').Value = '''0.999'
$ws.Range('E29').Value = '  -0.02%  '
$ws.Range('B30').Value = 'RenderToken'
$ws.Range('C30').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D30').Value = '''6.39'
$ws.Range('E30').Value = '  -0.83%  '
$ws.Range('D31').Value = '''7.25'
$ws.Range('E31').Value = '  +0.73%  '
$ws.Range('E32').Value = '  +3.72%  '
$ws.Range('B33').Value = 'EthereumClassic'
$ws.Range('C33').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D33').Value = '''20.92'
$ws.Range('E33').Value = '  +0.51%  '
$ws.Range('B34').Value = 'Fetch.AI'
$ws.Range('C34').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D34').Value = '''1.19'
$ws.Range('E34').Value = '  -0.51%  '
$ws.Range('D35').Value = '''158.19'
$ws.Range('E35').Value = '  +2.05%  '
$ws.Range('D36').Value = '''4.62'
$ws.Range('E36').Value = '  +1.36%  '
$ws.Range('D37').Value = '''6.15'
$ws.Range('E37').Value = '  +2.18%  '
$ws.Range('D38').Value = '''26.06'
$ws.Range('E38').Value = '  -4.74%  '
$ws.Range('D39').Value = '''1.25'
$ws.Range('E39').Value = '  -1.65%  '
$ws.Range('D40').Value = '''0.0664'
$ws.Range('E40').Value = '  -1.50%  '
$ws.Range('D41').Value = '''1.56'
$ws.Range('E41').Value = '  +12.91%  '
$ws.Range('D42').Value = '''4.00'
$ws.Range('E42').Value = '  +2.83%  '
$ws.Range('D43').Value = '''0.682'
$ws.Range('E43').Value = '  +4.07%  '
$ws.Range('D44').Value = '3.142.73'
$ws.Range('E44').Value = '  +0.83%  '
$ws.Range('D45').Value = '''36.80'
$ws.Range('E45').Value = '  +0.62%  '
$ws.Range('D46').Value = '''1.00'
$ws.Range('E46').Value = '  +0.06%  '
$ws.Range('B47').Value = 'VeChain'
$ws.Range('C47').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D47').Value = '''0.0262'
$ws.Range('E47').Value = '  +2.88%  '
$ws.Range('B48').Value = 'Maker'
$ws.Range('C48').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D48').Value = '2.277.32'
$ws.Range('E48').Value = '  -0.10%  '
$ws.Range('D49').Value = '''0.989'
$ws.Range('E49').Value = '  +5.15%  '
$ws.Range('B50').Value = 'Cosmos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D50').Value = '''6.07'
$ws.Range('E50').Value = '  +2.32%  '
$ws.Range('B51').Value = 'InjectiveProtocol'
$ws.Range('C51').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D51').Value = '''20.49'
$ws.Range('E51').Value = '  -0.37%  '
